$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet ("Export") is a flat list of account rows (Conta / Nome / Saldo).
# The update removes a handful of rows, inserts a couple of new ones, and
# relocates/updates the ADELE row with a new balance.
#
# All operations below are applied from the bottom of the sheet upward so
# that row numbers used for the earlier (still-untouched) rows stay valid.
# ---------------------------------------------------------------------------

# 1) Remove LARISSA / 004363260 (row 43) entirely.
$ws.Rows.Item(43).Delete()

# 2) Remove ALPHASITIO / 005305448 (row 17) entirely.
$ws.Rows.Item(17).Delete()

# 3) Replace the block MARCELO/ASSAKO/FABRICIO/ANDREA (rows 11-14) with a
#    single new row for MARIANA / 005046919 / 2800.
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "005046919"
$ws.Cells.Item(11, 2).Value = "MARIANA"
$ws.Cells.Item(11, 3).Value = 2800
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# 4) Remove ANA / 005198093 (row 8) entirely.
$ws.Rows.Item(8).Delete()

# 5) Remove the old ADELE / 004575632 row (row 7, balance 12353.49).
$ws.Rows.Item(7).Delete()

# 6) Insert the new ADELE / 004575632 row (balance 25354.8) right before
#    THOMAS, which is now row 5.
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004575632"
$ws.Cells.Item(5, 2).Value = "ADELE"
$ws.Cells.Item(5, 3).Value = 25354.8
